$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-21 Wednesday", "2024-02-22 Thursday"),
    @("12×26=", "51×27="),
    @("27×17=", "14×11="),
    @("73×31=", "99×83="),
    @("35×33=", "69×11="),
    @("79×46=", "75×52="),
    @("61×14=", "99×65="),
    @("51×70=", "51×24="),
    @("52×44=", "89×56="),
    @("41×81=", "88×40="),
    @("16×40=", "34×39="),
    @("51×55=", "34×98="),
    @("96×99=", "19×52="),
    @("63×40=", "55×93="),
    @("69×95=", "77×83="),
    @("26×54=", "34×11="),
    @("43×44=", "32×88="),
    @("19×90=", "98×97="),
    @("96×68=", "64×59="),
    @("52×88=", "52×34="),
    @("67×42=", "96×46="),
    @("12×93=", "96×93="),
    @("16×75=", "48×30="),
    @("57×85=", "39×53="),
    @("11×39=", "29×75="),
    @("43×32=", "63×81=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
